$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Capture the current (pre-insert) content of the VOUT "out(...)" row (row 50)
# before we overwrite it, so we can move it down into the newly inserted row.
$outCmd  = $ws.Range("B50").Value()
$outDesc = $ws.Range("C50").Value()

# Insert a new row below the "VOUT Command" header row (row 50), shifting
# everything from the old row 51 downward (xlShiftDown).
$ws.Rows("51:51").Insert(-4121)

# Row 50 now holds the new "range(c,r[,sn])" command entry.
$ws.Range("B50").Value = "range(c,r[,sn])"
$ws.Range("C50").Value = "for PHIDGET OUTPUT modules: sets voltage voltage range (r=5 fo r5V and r=10 for 10V)"

# The newly-inserted (blank) row 51 gets the command row that used to sit in
# row 50 ("out(<n>,<v>[,<sn>])").
$ws.Range("B51").Value = $outCmd
$ws.Range("C51").Value = $outDesc

# Restore selection state to match the post-edit workbook (selection now
# tracks the newly inserted row instead of the old B93:C94 leftover).
$ws.Range("B50:C50").Select()
$ws.Application.ActiveWindow.ScrollRow = 40

$wsButtons = $wb.Worksheets.Item("Buttons")
$wsButtons.Range("A1").Select()

$wsOptions = $wb.Worksheets.Item("Options")
$wsOptions.Range("B5").Select()

$ws.Activate()
$ws.Range("B50").Select()
